# "Added Category Parent ID"
# Fills in the Category ID / Parent ID pairs for the men's and women's
# "ready-to-wear-jeans" sub-categories (rows 564-589), removes the now
# superfluous trailing blank row (1000) and restores the view's
# selection to match where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- men-ready-to-wear-jeans children ---------------------------------
$ws.Range("A564").Value = "men-642"
$ws.Range("B564").Value = "men-ready-to-wear-jeans"
$ws.Range("A565").Value = "men-642-twin-pack"
$ws.Range("B565").Value = "men-ready-to-wear-jeans"
$ws.Range("A566").Value = "men-big-brother"
$ws.Range("B566").Value = "men-ready-to-wear-jeans"
$ws.Range("A567").Value = "men-big-jean"
$ws.Range("B567").Value = "men-ready-to-wear-jeans"
$ws.Range("A568").Value = "men-bob"
$ws.Range("B568").Value = "men-ready-to-wear-jeans"
$ws.Range("A569").Value = "men-bro"
$ws.Range("B569").Value = "men-ready-to-wear-jeans"
$ws.Range("A570").Value = "men-combat"
$ws.Range("B570").Value = "men-ready-to-wear-jeans"
$ws.Range("A571").Value = "men-cool-guy"
$ws.Range("B571").Value = "men-ready-to-wear-jeans"
$ws.Range("A572").Value = "men-eros"
$ws.Range("B572").Value = "men-ready-to-wear-jeans"
$ws.Range("A573").Value = "men-richard"
$ws.Range("B573").Value = "men-ready-to-wear-jeans"
$ws.Range("A574").Value = "men-sexy-twist"
$ws.Range("B574").Value = "men-ready-to-wear-jeans"
$ws.Range("A575").Value = "men-skater"
$ws.Range("B575").Value = "men-ready-to-wear-jeans"
$ws.Range("A576").Value = "men-super-twinky"
$ws.Range("B576").Value = "men-ready-to-wear-jeans"
$ws.Range("A577").Value = "men-tidy-biker"
$ws.Range("B577").Value = "men-ready-to-wear-jeans"

# --- women-ready-to-wear-jeans children --------------------------------
$ws.Range("A578").Value = "women-642"
$ws.Range("B578").Value = "women-ready-to-wear-jeans"
$ws.Range("A579").Value = "women-80s"
$ws.Range("B579").Value = "women-ready-to-wear-jeans"
$ws.Range("A580").Value = "women-boston"
$ws.Range("B580").Value = "women-ready-to-wear-jeans"
$ws.Range("A581").Value = "women-capri"
$ws.Range("B581").Value = "women-ready-to-wear-jeans"
$ws.Range("A582").Value = "women-cool-girl"
$ws.Range("B582").Value = "women-ready-to-wear-jeans"
$ws.Range("A583").Value = "women-eros"
$ws.Range("B583").Value = "women-ready-to-wear-jeans"
$ws.Range("A584").Value = "women-flare"
$ws.Range("B584").Value = "women-ready-to-wear-jeans"
$ws.Range("A585").Value = "women-jennifer"
$ws.Range("B585").Value = "women-ready-to-wear-jeans"
$ws.Range("A586").Value = "women-roadie"
$ws.Range("B586").Value = "women-ready-to-wear-jeans"
$ws.Range("A587").Value = "women-traveller"
$ws.Range("B587").Value = "women-ready-to-wear-jeans"
$ws.Range("A588").Value = "women--trumpet"
$ws.Range("B588").Value = "women-ready-to-wear-jeans"
$ws.Range("A589").Value = "women-twiggy"
$ws.Range("B589").Value = "women-ready-to-wear-jeans"

# --- drop the now-unused trailing blank row (shrinks A1:B1000 -> A1:B999)
$ws.Rows.Item(1000).Delete()

# --- restore view state (best-effort; scroll position is cosmetic) ----
$excel.ActiveWindow.ScrollRow = 559
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D585").Select()
